# Add second repository ("jiraflow-sample1") to the input sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the new row (row 3) with the second repository's data.
$ws.Range("A3").Value = "https://github.com/danielhantunes/jiraflow-sample1"
$ws.Range("B3").Value = "Joao"
$ws.Range("C3").Value = "email@example.com"

# Turn the repo URL and email cells into hyperlinks, same as row 2.
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/danielhantunes/jiraflow-sample1")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:email@example.com")

# Match the hyperlink cell style used by row 2 (applied after adding the
# hyperlinks, since Hyperlinks.Add restyles the cell on its own).
$ws.Range("A3").Style = $ws.Range("A2").Style
$ws.Range("C3").Style = $ws.Range("C2").Style

# Reflect the active selection shown in the saved workbook.
$ws.Range("C3").Select()
